$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $rng = $ws.Range($addr)
    $origStyle = $rng.Style
    $rng.NumberFormat = '@'
    $rng.Value = $val
    $rng.Style = $origStyle
}

Set-TextValue 'D2' '28.036.92'
$ws.Range('E2').Value = '  -2.17%  '
Set-TextValue 'D3' '1.830.19'
$ws.Range('E3').Value = '  -1.17%  '
$ws.Range('E4').Value = '  -0.21%  '
Set-TextValue 'D5' '326.32'
$ws.Range('E5').Value = '  -2.92%  '
$ws.Range('E6').Value = '  -0.17%  '
Set-TextValue 'D7' '0.4635'
$ws.Range('E7').Value = '  -0.42%  '
Set-TextValue 'D8' '0.3868'
$ws.Range('E8').Value = '  -1.23%  '
Set-TextValue 'D9' '0.07868'
$ws.Range('E9').Value = '  -0.27%  '
Set-TextValue 'D10' '0.9603'
$ws.Range('E10').Value = '  -2.34%  '
$ws.Range('B11').Value = 'WrappedEther'
$ws.Range('C11').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D11' '2.040.11'
$ws.Range('E11').Value = '  +10.61%  '
$ws.Range('B12').Value = 'Solana'
$ws.Range('C12').Value = 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol'
Set-TextValue 'D12' '21.89'
$ws.Range('E12').Value = '  -1.45%  '
Set-TextValue 'D13' '5.659'
$ws.Range('E13').Value = '  -3.30%  '
Set-TextValue 'D14' '6.894'
$ws.Range('E14').Value = '  -1.86%  '
Set-TextValue 'D15' '0.06736'
$ws.Range('E15').Value = '  -0.43%  '
$ws.Range('E16').Value = '  -0.26%  '
Set-TextValue 'D17' '86.89'
$ws.Range('E17').Value = '  -0.89%  '
Set-TextValue 'D18' '0.000009933'
$ws.Range('E18').Value = '  -1.58%  '
Set-TextValue 'D19' '16.62'
$ws.Range('E19').Value = '  -2.20%  '
$ws.Range('E20').Value = '  -0.09%  '
Set-TextValue 'D21' '28.049.96'
$ws.Range('E21').Value = '  -2.09%  '
Set-TextValue 'D22' '5.303'
$ws.Range('E22').Value = '  -2.02%  '
Set-TextValue 'D23' '10.98'
$ws.Range('E23').Value = '  -2.69%  '
Set-TextValue 'D24' '2.095'
$ws.Range('E24').Value = '  -1.46%  '
Set-TextValue 'D25' '2.083.36'
$ws.Range('E25').Value = '  +1.46%  '
Set-TextValue 'D26' '153.72'
$ws.Range('E26').Value = '  +0.08%  '
Set-TextValue 'D27' '19.17'
$ws.Range('E27').Value = '  -1.30%  '
Set-TextValue 'D28' '5.737'
$ws.Range('E28').Value = '  -8.62%  '
Set-TextValue 'D29' '1.975'
$ws.Range('E29').Value = '  -2.44%  '
Set-TextValue 'D30' '117.24'
$ws.Range('E30').Value = '  -0.49%  '
Set-TextValue 'D31' '0.9366'
$ws.Range('E31').Value = '  -4.43%  '
Set-TextValue 'D32' '0.09264'
$ws.Range('E32').Value = '  -1.95%  '
Set-TextValue 'D33' '5.293'
$ws.Range('E33').Value = '  -1.79%  '
Set-TextValue 'D34' '1.316'
$ws.Range('E34').Value = '  -2.82%  '
Set-TextValue 'D35' '3.317'
$ws.Range('E35').Value = '  -5.48%  '
Set-TextValue 'D36' '0.05867'
$ws.Range('E36').Value = '  -4.50%  '
Set-TextValue 'D37' '0.02143'
$ws.Range('E37').Value = '  -2.43%  '
Set-TextValue 'D38' '1.143'
$ws.Range('E38').Value = '  -1.03%  '
Set-TextValue 'D39' '7.766'
$ws.Range('E39').Value = '  +2.27%  '
Set-TextValue 'D40' '0.5591'
$ws.Range('E40').Value = '  -2.09%  '
Set-TextValue 'D41' '9.884'
$ws.Range('E41').Value = '  -2.15%  '
Set-TextValue 'D42' '0.1759'
$ws.Range('E42').Value = '  -1.43%  '
Set-TextValue 'D43' '1.205'
$ws.Range('E43').Value = '  -3.57%  '
Set-TextValue 'D44' '11.59'
$ws.Range('E44').Value = '  -1.65%  '
Set-TextValue 'D45' '0.5271'
$ws.Range('E45').Value = '  -2.22%  '
Set-TextValue 'D46' '0.07023'
$ws.Range('E46').Value = '  -1.55%  '
Set-TextValue 'D47' '2.157'
$ws.Range('E47').Value = '  -8.32%  '
Set-TextValue 'D48' '1.830'
$ws.Range('E48').Value = '  -4.19%  '
Set-TextValue 'D49' '112.94'
$ws.Range('E49').Value = '  -1.86%  '
Set-TextValue 'D50' '1.000'
$ws.Range('E50').Value = '  -0.23%  '
$ws.Range('E51').Value = '  +0.08%  '
